# CMPT_353_Project_Report.docx edits:
#  - proofing/wording touch-ups in the body text ("a html" -> "an html",
#    missing comma, run-on sentence split, "a SQL" -> "an SQL", and the
#    grammar-checker proofErr markers around "similar to"/"Similar to"
#    disappear once those passages are re-typed)
#  - a second author (Liam Neufeld / lwn282 / 11232603) added to the
#    title-page header, alongside the existing Henry Fang / hef052 / 11233914

$d = $word.ActiveDocument

function Replace-Body-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND (body): $old"
    }
}

function Replace-Header-Text($range, $old, $new) {
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND (header): $old"
    }
}

# 1) "Basically, a html client" -> "Basically, an html client"
Replace-Body-Text `
    "Basically, a html client and a MySQL database" `
    "Basically, an html client and a MySQL database"

# 2) "staff and customers and seeing" -> "staff and customers, and seeing"
Replace-Body-Text `
    "staff and customers and seeing if expected updates are performed." `
    "staff and customers, and seeing if expected updates are performed."

# 3) "...even Google it seemed..." -> "...even Google. It seemed..."
Replace-Body-Text `
    "Used by big names like Toyota, Adobe, and even Google it seemed like a solid technology" `
    "Used by big names like Toyota, Adobe, and even Google. It seemed like a solid technology"

# 4) "most similar to the tables" - re-typed so the gramStart/gramEnd
#    proofErr wrapper around "similar to" is dropped; wording unchanged
Replace-Body-Text `
    "ollections would be most similar to the tables that exist in SQL databases" `
    "ollections would be most similar to the tables that exist in SQL databases"

# 5) "Similar to defining the columns of a SQL database" -> drop the
#    proofErr wrapper around "Similar to" and fix "a SQL" -> "an SQL"
Replace-Body-Text `
    "which are JSON objects that define the structure and contents of a Collection. Similar to defining the columns of a SQL database." `
    "which are JSON objects that define the structure and contents of a Collection. Similar to defining the columns of an SQL database."

# 6/7) Header: add second author's initials and student number
$hdr = $d.Sections.Item(1).Headers.Item(1)

Replace-Header-Text $hdr.Range "hef052" "hef052, lwn282"
Replace-Header-Text $hdr.Range "11233914 " "11233914, 11232603 "

Write-Output "DONE"
